$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Correct the Portuguese translation in B4 ("then" -> "então")
$ws.Range("B4").Value = 'Os cenários possuem a cláusula "então"?'

# Mark every checklist item as completed
$ws.Range("C3:C9").Value = "OK"

# Column B needs to re-fit after the text change (it had bestFit applied)
$ws.Columns("B").AutoFit() | Out-Null

# Leave the selection where the author last left it
$ws.Range("F8").Select() | Out-Null
